$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data: date 2019-10-16 (serial 43754) and mileage 512
$ws.Range("A20").Value = 43754
$ws.Range("B20").Value = 512

# Move the active selection to B21, mirroring where the user clicked next
$ws.Range("B21").Select()
